$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85; existing rows 85-99 shift down to 86-100.
$ws.Rows.Item(85).Insert()

# Fill in the new row 85 with the values from the diff.
$ws.Cells.Item(85, 1).Value = 9
$ws.Cells.Item(85, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(85, 3).Value = "Metropolitana"
$ws.Cells.Item(85, 4).Value = 44748
$ws.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85, 5).Value = 13
$ws.Cells.Item(85, 6).Value = 100112005
$ws.Cells.Item(85, 7).Value = "Puerro"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 70
$ws.Cells.Item(85, 11).Value = 8000
$ws.Cells.Item(85, 12).Value = 8000
$ws.Cells.Item(85, 13).Value = 8000
$ws.Cells.Item(85, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(85, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(85, 16).Value = 400
$ws.Cells.Item(85, 17).Value = 20
$ws.Cells.Item(85, 18).Value = "Hortaliza"
